$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.451.72'
$ws.Range("E2").Value = '  -1.08%  '
$ws.Range("D3").Value = '1.959.08'
$ws.Range("E3").Value = '  -3.64%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''244.42'
$ws.Range("E5").Value = '  -1.65%  '
$ws.Range("D6").Value = '''0.617'
$ws.Range("E6").Value = '  -3.04%  '
$ws.Range("D7").Value = '''58.70'
$ws.Range("E7").Value = '  -6.50%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '''0.375'
$ws.Range("E9").Value = '  -3.46%  '
$ws.Range("D10").Value = '''55.70'
$ws.Range("E10").Value = '  -4.44%  '
$ws.Range("D11").Value = '''0.0844'
$ws.Range("E11").Value = '  +6.07%  '
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("D13").Value = '''0.839'
$ws.Range("E13").Value = '  -7.23%  '
$ws.Range("D14").Value = '''21.96'
$ws.Range("E14").Value = '  -5.06%  '
$ws.Range("D15").Value = '2.250.36'
$ws.Range("E15").Value = '  -3.67%  '
$ws.Range("D16").Value = '''13.58'
$ws.Range("E16").Value = '  -5.38%  '
$ws.Range("D17").Value = '''5.36'
$ws.Range("E17").Value = '  -3.00%  '
$ws.Range("D18").Value = '1.989.56'
$ws.Range("E18").Value = '  -2.29%  '
$ws.Range("D19").Value = '36.373.38'
$ws.Range("E19").Value = '  -1.19%  '
$ws.Range("D20").Value = '0.0₃0882'
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").Value = '''70.31'
$ws.Range("E21").Value = '  -2.49%  '
$ws.Range("D22").Value = '''231.26'
$ws.Range("E22").Value = '  -2.05%  '
$ws.Range("D23").Value = '''5.08'
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("E25").Value = '  +0.53%  '
$ws.Range("E26").Value = '  -2.07%  '
$ws.Range("D27").Value = '''9.55'
$ws.Range("E27").Value = '  -2.19%  '
$ws.Range("D28").Value = '''164.69'
$ws.Range("E28").Value = '  +3.37%  '
$ws.Range("D29").Value = '''19.68'
$ws.Range("E29").Value = '  -2.50%  '
$ws.Range("E30").Value = '  -12.41%  '
$ws.Range("E31").Value = '  -2.02%  '
$ws.Range("D32").Value = '''1.16'
$ws.Range("E32").Value = '  -1.14%  '
$ws.Range("D33").Value = '''4.75'
$ws.Range("E33").Value = '  -5.81%  '
$ws.Range("D34").Value = '''0.0640'
$ws.Range("E34").Value = '  +3.50%  '
$ws.Range("D35").Value = '''4.36'
$ws.Range("E35").Value = '  -3.04%  '
$ws.Range("E36").Value = '  -2.22%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("E38").Value = '  -1.66%  '
$ws.Range("E39").Value = '  -8.51%  '
$ws.Range("E40").Value = '  -6.46%  '
$ws.Range("D41").Value = '''0.0983'
$ws.Range("E41").Value = '  -0.78%  '
$ws.Range("E42").Value = '  -4.34%  '
$ws.Range("E43").Value = '  -3.19%  '
$ws.Range("E44").Value = '  -1.65%  '
$ws.Range("D45").Value = '''15.82'
$ws.Range("E45").Value = '  -7.21%  '
$ws.Range("E46").Value = '  -7.17%  '
$ws.Range("D47").Value = '''7.47'
$ws.Range("E47").Value = '  -3.22%  '
$ws.Range("D48").Value = '''89.22'
$ws.Range("E48").Value = '  -4.59%  '
$ws.Range("D49").Value = '1.347.99'
$ws.Range("E49").Value = '  -1.45%  '
$ws.Range("D50").Value = '''2.82'
$ws.Range("E50").Value = '  -3.11%  '
$ws.Range("D51").Value = '''45.32'
$ws.Range("E51").Value = '  -0.01%  '
